# Apply the edit described by the diff:
#  - Shared string "CDF(Z) " (with trailing space) used in column A (rows 2-29)
#    is replaced by "CDF(Z)" (trailing space removed).
#  - The saved selection in the sheet view moves from K31 to N11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the trailing-space text in column A for all data rows (2 through 29)
for ($r = 2; $r -le 29; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq "CDF(Z) ") {
        $cell.Value = "CDF(Z)"
    }
}

# Update the active selection to N11, as recorded in the saved sheet view
$ws.Range("N11").Select()
